# Update the "Förändrad" (changed) date column (C2:C9) by one day:
# 45207 -> 45208 (2023-10-08 -> 2023-10-09)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C9").Value = 45208
